# booking.xlsx — "update check in to allow many rooms"
# Adds a new booking row (row 15) to Sheet1 and drops the stray per-cell
# styling that had crept onto a handful of E/K cells in rows 2-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new booking as row 15 -------------------------------------
# Seed row 15 from row 14 first so the date cells (C/D/G) inherit the exact
# same cell style (s="1") instead of minting a brand-new style entry.
$ws.Range("A14:K14").Copy($ws.Range("A15:K15")) | Out-Null

$ws.Range("A15").Value2 = 14
$ws.Range("B15").Value2 = 320
$ws.Range("C15").Value2 = 44930
$ws.Range("D15").Value2 = 44932
$ws.Range("E15").Formula = "=(20*B15)/100"
$ws.Range("F15").Value2 = 64
$ws.Range("G15").Value2 = 44916
$ws.Range("H15").Value2 = $ws.Range("H11").Value2
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 4
$ws.Range("K15").Value2 = 107275663

# --- Strip the leftover formatting on the affected E/K cells --------------
# These cells carried a redundant explicit style (identical to the default)
# which gets cleaned up as part of this edit. (ClearFormats only honours the
# first area of a multi-area Range here, so clear each cell individually.)
foreach ($addr in @("E2", "K2", "E4", "E5", "E6", "E7", "K7", "K8")) {
    $ws.Range($addr).ClearFormats() | Out-Null
}

# --- Update the current selection to the newly added row ------------------
$ws.Range("M15").Select() | Out-Null
